$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two data rows (rows 7 and 8), which shifts the remaining
# rows up by two.
$ws.Rows("7:8").Delete()

# Renumber the "م" (sequence number) column for the rows that moved up.
$ws.Range("A7").Value = 1
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3

# Update the trailing total (sum of the "سعر البيع" column) to reflect the
# removal of the two deleted rows' sale prices (218 + 220).
$ws.Range("P10").Value = 2235

# Match the total row's final height.
$ws.Rows("10").RowHeight = 24.75
